# TC18_Canine_Filter_Breed-Chinese.xlsx - corrected ICDC Breed 1-14 scripts
#
# The FilesTab query (row 4, column B) is corrected: the `File Type` and
# `Breed` columns are dropped from the RETURN clause (and their matching
# WITH/coalesce lines removed), matching the same cleanup already applied
# to the other tabs' scripts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFilesQuery = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN [''Chinese Shar-Pei'']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '''') AS `File Name`,
          coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(f.file_format, '''') AS `Format`,
        coalesce(f.file_size, '''') AS `Size`,
        coalesce(c.case_id, '''') AS `Case ID`,
         coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'''') AS `Study Code`'

# Update the FilesTab row's query cell (B4) with the corrected script.
$ws.Range("B4").Value2 = $newFilesQuery

# The row auto-shrinks slightly now that two lines were removed from the
# wrapped text.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection / scroll moved down to the FilesTab row after the edit.
[void]$ws.Range("B4").Select()
